# Update GMM (column G) values in Sheet1 per the source data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value  = 2974.769829392405
$ws.Range("G4").Value  = 10892.8757319292
$ws.Range("G5").Value  = 14414.62280778821
$ws.Range("G6").Value  = 0.4784345194989783
$ws.Range("G8").Value  = 0.5207484758709566
$ws.Range("G9").Value  = 0.5604348830544731
$ws.Range("G10").Value = 0.4157836489644302
$ws.Range("G12").Value = 0.4063563102698645
$ws.Range("G13").Value = 0.5444396508305743
